$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "County" column between Address (B) and Program Name (old C) ---
$ws.Columns("C:C").Insert()

# Header
$ws.Range("C1").Value = "County"

# County per row - all schools are in Orange County except Longfellow Elementary,
# whose address (3610 Eucalyptus Ave, Riverside, CA 92507) is in Riverside.
$ws.Range("C2").Value = "Orange County"
$ws.Range("C3").Value = "Orange County"
$ws.Range("C4").Value = "Orange County"
$ws.Range("C5").Value = "Orange County"
$ws.Range("C6").Value = "Orange County"
$ws.Range("C7").Value = "Orange County"
$ws.Range("C8").Value = "Orange County"
$ws.Range("C9").Value = "Orange County"
$ws.Range("C10").Value = "Orange County"
$ws.Range("C11").Value = "Orange County"
$ws.Range("C12").Value = "Riverside"

# Restore the old column-3/4 formatting width (20) over the new County column, matching
# the width that already spanned that region before the insert.
$ws.Columns("C:C").ColumnWidth = 19.17

# Widen column B (Address) and auto-size the new last column (Hours moved to G)
$ws.Columns("B:B").ColumnWidth = 41.665
$ws.Columns("G:G").ColumnWidth = 18.83

# --- View state: drop the old horizontal scroll, zoom to 79%, move selection to C14 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$win.Zoom = 79

$ws.Range("C14").Select() | Out-Null
